{"js": "// Replace \"Paradigm\" with \"AspenTech\" in the acknowledgements paragraph\n// (the sponsor being thanked for the SKUA-Gocad software / dev kit).\nconst body = context.document.body;\nconst results = body.search(\"Paradigm\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"AspenTech\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace \"Paradigm\" with \"AspenTech\" in the acknowledgements paragraph\n# (the sponsor being thanked for the SKUA-Gocad software / dev kit).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Paradigm\"\n$find.Replacement.Text = \"AspenTech\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceAll) | Out-Null\n"}
